{"js": "// Each table cell holds one \"a\u00f7b=c, d\" style division fact. The commit\n// regenerates the worksheet, so every fact is swapped for a freshly\n// generated one. Match (and replace) each old fact text with its new\n// counterpart via Word.js' body.search().\nconst body = context.document.body;\nconst replacements = [\n  [\"325\u00f75=65, 0\", \"454\u00f75=90, 4\"],\n  [\"377\u00f79=41, 8\", \"582\u00f72=291, 0\"],\n  [\"814\u00f73=271, 1\", \"814\u00f74=203, 2\"],\n  [\"476\u00f72=238, 0\", \"314\u00f75=62, 4\"],\n  [\"410\u00f72=205, 0\", \"607\u00f76=101, 1\"],\n  [\"511\u00f79=56, 7\", \"476\u00f79=52, 8\"],\n  [\"810\u00f72=405, 0\", \"198\u00f74=49, 2\"],\n  [\"225\u00f75=45, 0\", \"598\u00f79=66, 4\"],\n  [\"990\u00f73=330, 0\", \"555\u00f73=185, 0\"],\n  [\"951\u00f74=237, 3\", \"140\u00f78=17, 4\"],\n  [\"175\u00f75=35, 0\", \"597\u00f78=74, 5\"],\n  [\"844\u00f79=93, 7\", \"672\u00f79=74, 6\"],\n  [\"529\u00f75=105, 4\", \"173\u00f77=24, 5\"],\n  [\"276\u00f72=138, 0\", \"826\u00f76=137, 4\"],\n  [\"701\u00f78=87, 5\", \"542\u00f76=90, 2\"],\n  [\"286\u00f79=31, 7\", \"830\u00f75=166, 0\"],\n  [\"216\u00f79=24, 0\", \"402\u00f79=44, 6\"],\n  [\"794\u00f73=264, 2\", \"596\u00f74=149, 0\"],\n  [\"425\u00f79=47, 2\", \"708\u00f73=236, 0\"],\n  [\"712\u00f74=178, 0\", \"876\u00f77=125, 1\"],\n  [\"508\u00f74=127, 0\", \"492\u00f76=82, 0\"],\n  [\"638\u00f77=91, 1\", \"823\u00f73=274, 1\"],\n  [\"452\u00f73=150, 2\", \"529\u00f79=58, 7\"],\n  [\"454\u00f79=50, 4\", \"705\u00f75=141, 0\"],\n  [\"110\u00f78=13, 6\", \"795\u00f73=265, 0\"],\n];\n\n// All old/new fact strings are mutually distinct substrings of one\n// another (verified offline), so a simple sequential search/replace\n// loop cannot cross-contaminate an earlier or later substitution.\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Each table cell holds one \"a\u00f7b=c, d\" style division fact. The commit\n# regenerates the worksheet, so every fact is swapped for a freshly\n# generated one. Replace each old fact text with its new counterpart via\n# Word's Find/Replace (wdReplaceOne), scanning the whole document body.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('325\u00f75=65, 0', '454\u00f75=90, 4'),\n    @('377\u00f79=41, 8', '582\u00f72=291, 0'),\n    @('814\u00f73=271, 1', '814\u00f74=203, 2'),\n    @('476\u00f72=238, 0', '314\u00f75=62, 4'),\n    @('410\u00f72=205, 0', '607\u00f76=101, 1'),\n    @('511\u00f79=56, 7', '476\u00f79=52, 8'),\n    @('810\u00f72=405, 0', '198\u00f74=49, 2'),\n    @('225\u00f75=45, 0', '598\u00f79=66, 4'),\n    @('990\u00f73=330, 0', '555\u00f73=185, 0'),\n    @('951\u00f74=237, 3', '140\u00f78=17, 4'),\n    @('175\u00f75=35, 0', '597\u00f78=74, 5'),\n    @('844\u00f79=93, 7', '672\u00f79=74, 6'),\n    @('529\u00f75=105, 4', '173\u00f77=24, 5'),\n    @('276\u00f72=138, 0', '826\u00f76=137, 4'),\n    @('701\u00f78=87, 5', '542\u00f76=90, 2'),\n    @('286\u00f79=31, 7', '830\u00f75=166, 0'),\n    @('216\u00f79=24, 0', '402\u00f79=44, 6'),\n    @('794\u00f73=264, 2', '596\u00f74=149, 0'),\n    @('425\u00f79=47, 2', '708\u00f73=236, 0'),\n    @('712\u00f74=178, 0', '876\u00f77=125, 1'),\n    @('508\u00f74=127, 0', '492\u00f76=82, 0'),\n    @('638\u00f77=91, 1', '823\u00f73=274, 1'),\n    @('452\u00f73=150, 2', '529\u00f79=58, 7'),\n    @('454\u00f79=50, 4', '705\u00f75=141, 0'),\n    @('110\u00f78=13, 6', '795\u00f73=265, 0')\n)\n\n# wdReplaceOne = 1, wdFindWrapContinue = 1 (Wrap), wdReplace = 2 (here\n# passed positionally below). All old/new fact strings are mutually\n# distinct substrings of one another (verified offline), so running the\n# replacements sequentially over $d.Content cannot cross-contaminate an\n# earlier or later substitution.\n$wdReplaceOne = 1\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $found = $range.Find.Execute(\n        $oldText,   # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $newText,   # ReplaceWith\n        $wdReplaceOne\n    )\n\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
